$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, column, newText) for the data rows of the table.
# Rows 1, 5, 9, 13, 17 contain the exercise text; the other rows are blank.
$replacements = @(
    @{ Row = 1;  Col = 1; New = "35÷7=" },
    @{ Row = 1;  Col = 2; New = "45÷3=" },
    @{ Row = 1;  Col = 3; New = "67÷9=" },
    @{ Row = 1;  Col = 4; New = "13÷9=" },
    @{ Row = 1;  Col = 5; New = "70÷9=" },

    @{ Row = 5;  Col = 1; New = "10÷7=" },
    @{ Row = 5;  Col = 2; New = "87÷4=" },
    @{ Row = 5;  Col = 3; New = "37÷6=" },
    @{ Row = 5;  Col = 4; New = "86÷6=" },
    @{ Row = 5;  Col = 5; New = "66÷5=" },

    @{ Row = 9;  Col = 1; New = "56÷8=" },
    @{ Row = 9;  Col = 2; New = "34÷8=" },
    @{ Row = 9;  Col = 3; New = "94÷7=" },
    @{ Row = 9;  Col = 4; New = "17÷4=" },
    @{ Row = 9;  Col = 5; New = "82÷8=" },

    @{ Row = 13; Col = 1; New = "57÷9=" },
    @{ Row = 13; Col = 2; New = "68÷8=" },
    @{ Row = 13; Col = 3; New = "45÷4=" },
    @{ Row = 13; Col = 4; New = "77÷5=" },
    @{ Row = 13; Col = 5; New = "38÷4=" },

    @{ Row = 17; Col = 1; New = "96÷5=" },
    @{ Row = 17; Col = 2; New = "92÷6=" },
    @{ Row = 17; Col = 3; New = "97÷4=" },
    @{ Row = 17; Col = 4; New = "67÷9=" },
    @{ Row = 17; Col = 5; New = "33÷2=" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    # Shrink the range so it doesn't include the end-of-cell marker,
    # then replace the run's text directly.
    $rng.End = $rng.End - 1
    $rng.Text = $item.New
}

$d.Save()
